{"js": "// Apply the \"Add feria in DRE, MER and ML\" field-reordering edits to the\n// logical-model table descriptions in the document body.\n//\n// Each edit below targets a short, unique substring within one paragraph's\n// parenthesized attribute list and replaces it with the reordered /\n// renamed text described by the diff.\n\nconst replacements = [\n  // usuario(id, email, senha, nome, cpf) -> usuario(id, email, password, nome, cpf)\n  { find: \", senha, nome, \", replace: \", password, nome, \" },\n\n  // cabeleireiro(id , cargo n\u00famero, rua, cidade, estado)\n  //   -> cabeleireiro(id , rua, n\u00famero, cidade, estado)\n  { find: \" , cargo n\u00famero, rua, cidade, estado)\", replace: \" , rua, n\u00famero, cidade, estado)\" },\n\n  // atendimento(id, horario, status, data, id_servico, id_cliente, id_cabeleireiro)\n  //   -> atendimento(id, horario, data, id_servico, id_cliente, id_cabeleireiro, status)\n  { find: \", status, data, \", replace: \", data, \" },\n  { find: \"id_cabeleireiro) \", replace: \"id_cabeleireiro, status) \" },\n\n  // comentario(id, descricao, data, titulo, id_cliente, id_servico)\n  //   -> comentario(id, titulo, descricao, data, id_cliente, id_servico)\n  { find: \"descricao, data, titulo, \", replace: \"titulo, descricao, data, \" },\n\n  // produto(id, marca, preco, validade, quantidade)\n  //   -> produto(id, nome, marca, preco, validade, quantidade)\n  // (search starts right after the bold/underlined \"id\" run so that run's\n  // formatting is left untouched)\n  { find: \", marca, \", replace: \", nome, marca, \" },\n\n  // fornecedor(id, email, telefone, nome, cnpj, cpf, id_produto)\n  //   -> fornecedor(id, nome, email, telefone, cnpj, cpf, id_produto)\n  // (same reasoning: start after the bold/underlined \"id\" run)\n  { find: \", email, telefone, nome, cnpj\", replace: \", nome, email, telefone, cnpj\" },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${JSON.stringify(find)}`);\n  }\n\n  // Each search string is unique in the document, but replace every match\n  // defensively in case of repeats.\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Add feria in DRE, MER and ML\" field-reordering edits to the\n# logical-model table descriptions in the document body.\n#\n# Each block below performs a plain text Find/Replace of a short, unique\n# substring within one paragraph's parenthesized attribute list, swapping\n# in the reordered / renamed text described by the diff.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $found = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Search text not found: $findText\"\n    }\n}\n\n# usuario(id, email, senha, nome, cpf) -> usuario(id, email, password, nome, cpf)\nReplace-Text \", senha, nome, \" \", password, nome, \"\n\n# cabeleireiro(id , cargo n\u00famero, rua, cidade, estado)\n#   -> cabeleireiro(id , rua, n\u00famero, cidade, estado)\nReplace-Text \" , cargo n\u00famero, rua, cidade, estado)\" \" , rua, n\u00famero, cidade, estado)\"\n\n# atendimento(id, horario, status, data, id_servico, id_cliente, id_cabeleireiro)\n#   -> atendimento(id, horario, data, id_servico, id_cliente, id_cabeleireiro, status)\nReplace-Text \", status, data, \" \", data, \"\nReplace-Text \"id_cabeleireiro) \" \"id_cabeleireiro, status) \"\n\n# comentario(id, descricao, data, titulo, id_cliente, id_servico)\n#   -> comentario(id, titulo, descricao, data, id_cliente, id_servico)\nReplace-Text \"descricao, data, titulo, \" \"titulo, descricao, data, \"\n\n# produto(id, marca, preco, validade, quantidade)\n#   -> produto(id, nome, marca, preco, validade, quantidade)\n# (search starts right after the bold/underlined \"id\" run so that run's\n# formatting is left untouched)\nReplace-Text \", marca, \" \", nome, marca, \"\n\n# fornecedor(id, email, telefone, nome, cnpj, cpf, id_produto)\n#   -> fornecedor(id, nome, email, telefone, cnpj, cpf, id_produto)\n# (same reasoning: start after the bold/underlined \"id\" run)\nReplace-Text \", email, telefone, nome, cnpj\" \", nome, email, telefone, cnpj\"\n"}
